$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.792.55"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.637.15"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.49"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  -0.73%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0791"
$ws.Range("E11").Value = "  +1.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.28"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "1.862.79"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "1.636.91"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("E16").Value = "  -0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.18"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "25.819.44"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.13"
$ws.Range("E21").Value = "  -0.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.98"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.40"
$ws.Range("E23").Value = "  +2.63%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  +2.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.53"
$ws.Range("E26").Value = "  +3.13%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.55"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.906"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").Value = "1.133.68"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.43"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("D45").Value = "1.772.34"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("E48").Value = "  -1.61%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.41"
$ws.Range("E50").Value = "  +2.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.49"
$ws.Range("E51").Value = "  -3.26%  "
